# Applies the "NIT-9009959885" EC update:
#  - Adds a new worker (ALEJANDRO MANUEL ACOSTA PACHECO) with 3 mora periods
#    at the top of the detail table (rows 16-18).
#  - Keeps the existing worker (NEIVER JOSE ALVAREZ PAYARES) but re-orders his
#    52 mora periods from ascending to descending (most recent period first).
#  - Updates the summary fields (Valor Mora total, Cant. Trabajadores).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert 3 new rows right above the existing first detail row (row 16),
#    pushing the existing worker's 52 rows + the signature block down by 3.
# ---------------------------------------------------------------------------
$ws.Rows("16:18").Insert()

# The 3 freshly-inserted rows inherit the header row's formatting; copy the
# formatting from an existing detail row (now at row 19) onto them instead.
$ws.Rows("19").Copy()
$ws.Rows("16:18").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Fill in the new worker's 3 rows (16-18): ALEJANDRO MANUEL ACOSTA PACHECO
# ---------------------------------------------------------------------------
$newRows = @(
    @{ Row = 16; Period = "2310"; Mora = 35574; Salario = 1160000 },
    @{ Row = 17; Period = "2309"; Mora = 46400; Salario = 1160000 },
    @{ Row = 18; Period = "2308"; Mora = 46400; Salario = 1160000 }
)
foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 2).Value = "CC"
    $ws.Cells.Item($r.Row, 3).Value = "73125996"
    $ws.Cells.Item($r.Row, 4).Value = "ALEJANDRO MANUEL ACOSTA PACHECO"
    $ws.Cells.Item($r.Row, 5).Value = $r.Period
    $ws.Cells.Item($r.Row, 6).Value = $r.Mora
    $ws.Cells.Item($r.Row, 7).Value = $r.Salario
}

# ---------------------------------------------------------------------------
# 3. Re-order the existing worker's periods (now rows 19-70) from ascending
#    (1907 -> 2310) to descending (2310 -> 1907). Only the Period (E) and
#    Valor Mora (F) columns change value/order; B/C/D/G stay identical for
#    every row, since it's still the same worker with the same salary.
# ---------------------------------------------------------------------------
$periods = @()
foreach ($n in 10..1)  { $periods += ("23{0:D2}" -f $n) }
foreach ($n in 12..1)  { $periods += ("22{0:D2}" -f $n) }
foreach ($n in 12..1)  { $periods += ("21{0:D2}" -f $n) }
foreach ($n in 12..1)  { $periods += ("20{0:D2}" -f $n) }
foreach ($n in 12..7)  { $periods += ("19{0:D2}" -f $n) }

for ($i = 0; $i -lt $periods.Count; $i++) {
    $row = 19 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
    if ($i -eq 0) {
        $ws.Cells.Item($row, 6).Value = 25396
    } else {
        $ws.Cells.Item($row, 6).Value = 33125
    }
}

# ---------------------------------------------------------------------------
# 4. Update the summary block.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 1843145   # VALOR MORA total
$ws.Range("C13").Value = 2         # Cant. Trabajadores (was 1, now 2)
# F13 (Cant. Periodos = 52) is unchanged.
